$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36 (STOP_HEAD row): replace the old "Entschuldigung!" / "Sorry!"
#     heading with the new participant greeting. ---
$ws.Range("B36").Value = 'Liebe:r Teilnehmer:in,'
$ws.Range("C36").Value = 'Dear participant,'

# --- Row 37 (STOP_TEXT row): replace the old "listening device is not
#     suitable" message with the new, longer stop-page copy. Column C
#     (English) is written before column B (German) so the workbook's
#     shared-string table ends up in the same order as produced by Excel. ---
$ws.Range("C37").Value = 'the interview has ended. This can have several reasons:<ol><li>Sufficient test persons with a similar playback device have already participated.</li><li> Your playback device is not part of our target group.</li></ol>Since future studies may require more subjects and different target groups, we would be happy to see you again in the next study.\\Thank you for your interest and participation.
'
$ws.Range("B37").Value = 'die Befragung wurde beendet. Das kann mehrere Gründe haben:<ol><li>Es haben bereits genügend Probanden mit einem ähnlichen Wiedergabegerät teilgenommen.</li><li> Ihr Wiedergabegerät ist nicht Teil unserer Zielgruppe.</li>
</ol>Da zukünftige Untersuchungen vielleicht mehr Probanden und andere Zielgruppen erfordern, würden wir uns freuen, Sie in der nächsten Studie wieder begrüßen zu können. \\Wir bedanken uns für Ihr Interesse und Ihre Teilnahme.'

# The new STOP_TEXT copy is much longer, so wrap it and size the row the
# same way the other long-text rows on this sheet are sized.
$ws.Range("B37:C37").WrapText = $true
$ws.Rows(37).RowHeight = 180
